$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("ECs", "Adm", "Ramp3", "ECs", 3, 1, 12.31940633333333, 36.958219, 0.3801768364207869, 0.3801768364207869, 2, 0.6666666666666666, 0.4857986666666667, 1.457396, 0.4938122518903432, 0.4938122518903432, 5.984751170858223, 53.86276053772401, 0.1877359797094954, 0.1877359797094954)
    ,@("ECs", "Adm", "Ramp3", "FAPs", 3, 1, 12.31940633333333, 36.958219, 0.3801768364207869, 0.3801768364207869, 3, 1, 0.4477856666666667, 1.343357, 0.4551722011468782, 0.4551722011468782, 5.516453577909222, 49.648082201183, 0.1730459274587062, 0.1730459274587062)
    ,@("ECs", "Adm", "Ramp3", "MuSCs", 3, 1, 12.31940633333333, 36.958219, 0.3801768364207869, 0.3801768364207869, 1, 0.3333333333333333, 0.05018766666666667, 0.150563, 0.05101554696277864, 0.05101554696277864, 0.6182822585885556, 5.564540327297, 0.01939492925258526, 0.01939492925258526)
    ,@("FAPs", "Adm", "Ramp3", "ECs", 3, 1, 18.11265066666666, 54.33795199999999, 0.5589563363143816, 0.5589563363143816, 2, 0.6666666666666666, 0.4857986666666667, 1.457396, 0.4938122518903432, 0.4938122518903432, 8.799101543665778, 79.19191389299201, 0.2760194871437808, 0.2760194871437808)
    ,@("FAPs", "Adm", "Ramp3", "FAPs", 3, 1, 18.11265066666666, 54.33795199999999, 0.5589563363143816, 0.5589563363143816, 3, 1, 0.4477856666666667, 1.343357, 0.4551722011468782, 0.4551722011468782, 8.110585353873777, 72.995268184864, 0.2544213859452118, 0.2544213859452118)
    ,@("FAPs", "Adm", "Ramp3", "MuSCs", 3, 1, 18.11265066666666, 54.33795199999999, 0.5589563363143816, 0.5589563363143816, 1, 0.3333333333333333, 0.05018766666666667, 0.150563, 0.05101554696277864, 0.05101554696277864, 0.9090316741084442, 8.181285066975999, 0.02851546322538902, 0.02851546322538902)
    ,@("MuSCs", "Adm", "Ramp3", "ECs", 3, 1, 1.603212, 4.809636, 0.04947511672073613, 0.04947511672073613, 2, 0.6666666666666666, 0.4857986666666667, 1.457396, 0.4938122518903432, 0.4938122518903432, 0.7788382519840001, 7.009544267856001, 0.02443141880040428, 0.02443141880040428)
    ,@("MuSCs", "Adm", "Ramp3", "FAPs", 3, 1, 1.603212, 4.809636, 0.04947511672073613, 0.04947511672073613, 3, 1, 0.4477856666666667, 1.343357, 0.4551722011468782, 0.4551722011468782, 0.7178953542280001, 6.461058188052001, 0.02251969777977618, 0.02251969777977618)
    ,@("MuSCs", "Adm", "Ramp3", "MuSCs", 3, 1, 1.603212, 4.809636, 0.04947511672073613, 0.04947511672073613, 1, 0.3333333333333333, 0.05018766666666667, 0.150563, 0.05101554696277864, 0.05101554696277864, 0.080461469452, 0.7241532250680001, 0.002524000140555668, 0.002524000140555668)
    ,@("Resolving-Mac", "Adm", "Ramp3", "ECs", 3, 1, 0.3691416666666666, 1.107425, 0.01139171054409548, 0.01139171054409548, 2, 0.6666666666666666, 0.4857986666666667, 1.457396, 0.4938122518903432, 0.4938122518903432, 0.1793285294777778, 1.6139567653, 0.005625366236662756, 0.005625366236662754)
    ,@("Resolving-Mac", "Adm", "Ramp3", "FAPs", 3, 1, 0.3691416666666666, 1.107425, 0.01139171054409548, 0.01139171054409548, 3, 1, 0.4477856666666667, 1.343357, 0.4551722011468782, 0.4551722011468782, 0.1652963473027778, 1.487667125725, 0.005185189963184041, 0.005185189963184041)
    ,@("Resolving-Mac", "Adm", "Ramp3", "MuSCs", 3, 1, 0.3691416666666666, 1.107425, 0.01139171054409548, 0.01139171054409548, 1, 0.3333333333333333, 0.05018766666666667, 0.150563, 0.05101554696277864, 0.05101554696277864, 0.01852635891944444, 0.166737230275, 0.0005811543442486836, 0.0005811543442486834)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowVals[$c]
    }
}